$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price / Volume) to Text format so that
# numeric-looking values (e.g. "1.00", "60.360.23") are stored as
# literal strings, matching the source data which is all inlineStr text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "60.360.23"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "2.887.63"
$ws.Range("E3").Value = "  -3.36%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "526.66"
$ws.Range("E5").Value = "  -4.37%  "
$ws.Range("D6").Value = "141.97"
$ws.Range("E6").Value = "  -6.27%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").Value = "2.889.22"
$ws.Range("E9").Value = "  -3.18%  "
$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  -4.04%  "
$ws.Range("D11").Value = "5.90"
$ws.Range("E11").Value = "  -5.11%  "
$ws.Range("D12").Value = "0.357"
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("D13").Value = "3.382.39"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").Value = "60.335.85"
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("D16").Value = "22.51"
$ws.Range("E16").Value = "  -4.65%  "
$ws.Range("D17").Value = "2.873.11"
$ws.Range("E17").Value = "  -4.11%  "
$ws.Range("D18").Value = "0.0000141"
$ws.Range("E18").Value = "  -4.53%  "
$ws.Range("D19").Value = "4.93"
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").Value = "11.54"
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("D21").Value = "360.68"
$ws.Range("E21").Value = "  -7.52%  "
$ws.Range("D22").Value = "6.53"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "63.26"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").Value = "3.007.66"
$ws.Range("E25").Value = "  -4.11%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "0.449"
$ws.Range("E26").Value = "  -3.45%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.181"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "7.77"
$ws.Range("E29").Value = "  -7.85%  "
$ws.Range("D30").Value = "0.0₃0853"
$ws.Range("E30").Value = "  -9.67%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "1.67"
$ws.Range("E32").Value = "  -3.54%  "
$ws.Range("D33").Value = "19.44"
$ws.Range("E33").Value = "  -4.68%  "
$ws.Range("D34").Value = "149.69"
$ws.Range("E34").Value = "  -5.68%  "
$ws.Range("D35").Value = "4.33"
$ws.Range("E35").Value = "  -6.73%  "
$ws.Range("D36").Value = "5.55"
$ws.Range("E36").Value = "  -7.12%  "
$ws.Range("D37").Value = "0.988"
$ws.Range("E37").Value = "  -7.37%  "
$ws.Range("D38").Value = "1.20"
$ws.Range("E38").Value = "  -6.26%  "
$ws.Range("D39").Value = "37.76"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "1.49"
$ws.Range("E40").Value = "  -4.34%  "
$ws.Range("D41").Value = "2.320.65"
$ws.Range("E41").Value = "  -4.94%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "3.65"
$ws.Range("E42").Value = "  -6.14%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.640"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").Value = "20.59"
$ws.Range("E44").Value = "  -7.30%  "
$ws.Range("D45").Value = "0.0571"
$ws.Range("E45").Value = "  -3.51%  "
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "5.03"
$ws.Range("E47").Value = "  +2.93%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0233"
$ws.Range("E48").Value = "  -4.63%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "10.34"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").Value = "0.0927"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").Value = "248.98"
$ws.Range("E51").Value = "  -4.63%  "

# Restore default (Normal) style so no stray number-format style lingers
# on cells beyond what the original workbook had.
$ws.Range("D2:E51").Style = "Normal"
